$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of test-case data (ids 22..29 in column A already present).
# Columns: B = element name, C = case description.
$rows = @(
    @{ Row = 23; B = "CircleParallel"; C = "circle inside fov" },
    @{ Row = 24; B = "CircleParallel"; C = "circle cut by fov" },
    @{ Row = 25; B = "DialDeg";        C = "dial rail on CircleParallel inside fov" },
    @{ Row = 26; B = "DialDeg";        C = "dial rail on CircleParallel cut by fov" },
    @{ Row = 27; B = "CircleMeridian"; C = "circle inside fov" },
    @{ Row = 28; B = "CircleMeridian"; C = "circle cut by fov" },
    @{ Row = 29; B = "DialDeg";        C = "dial rail on CircleMeridian inside fov" },
    @{ Row = 30; B = "DialDeg";        C = "dial rail on CircleMeridian cut by fov" }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Range("B$r").Value = $item.B
    $ws.Range("C$r").Value = $item.C
    $ws.Range("E$r").Formula = '="tf-"&TEXT(A' + $r + ',"0000")&".xml"'
}

# Selection moved to E29 as part of this edit.
[void]$ws.Range("E29").Select()
